# Remove the trailing "Ver no Jupiter..." and "(c) 2020 ..." paragraphs
# (site-footer boilerplate removed from the generated page) while leaving
# everything else — including the blank paragraph that precedes them and
# the page-break paragraph that follows — untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$find1 = $d.Content
$found1 = $find1.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found1) {
    throw "Could not find the 'Ver no Jupiter...' paragraph"
}
$firstIndex = $find1.Paragraphs.Item(1).Index

# Locate the "(c) 2020 . Contact: ..." paragraph.
$find2 = $d.Content
$found2 = $find2.Find.Execute(
    "Contact: luizeleno@usp.br",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found2) {
    throw "Could not find the copyright footer paragraph"
}
$lastIndex = $find2.Paragraphs.Item(1).Index

# Delete both paragraphs (including their paragraph marks) in one go.
$startPara = $d.Paragraphs.Item($firstIndex)
$endPara = $d.Paragraphs.Item($lastIndex)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
